# Add files via upload
#
# 1) Fill in the "consolidated" sheet (sheet1) with the rolled-up estimates.
# 2) Add a new "Ruoxuan Li" worksheet (sheet3) with her own estimates in
#    columns A (story) and D (estimate).
# Sheet2 ("Fanjie Gao") keeps the same visible content; it is not touched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: "consolidated"
# ---------------------------------------------------------------------
$consolidated = $wb.Worksheets.Item(1)
$fanjie = $wb.Worksheets.Item(2)

$consolidatedRows = @(
    @("User Story 1", "time estimate"),
    @("a frontend to place the app", "3h"),
    @("OAuth to spotify", "2.5h"),
    @("store the token", "2h"),
    @("test api for spotify", "2.5h"),
    @("test api for imdb", "2.5h"),
    @("get user playlist genre", "0.5h"),
    @("query imdb with the genres", "0.5h"),
    @("algorithm to sort the movies", "1h"),
    @("database to store user and their movie list", "3h"),
    @("option to remove movie form list", "1h"),
    @("excluded movie list in database", "1h")
)

# Copy the existing cell formatting (font/style) from the already-populated
# "Fanjie Gao" sheet so the new cells use the same style index instead of
# creating brand new style entries.
$fanjie.Range("A1:B12").Copy() | Out-Null
$consolidated.Range("A1:B12").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $consolidatedRows.Count; $i++) {
    $r = $i + 1
    $consolidated.Cells.Item($r, 1).Value = $consolidatedRows[$i][0]
    $consolidated.Cells.Item($r, 2).Value = $consolidatedRows[$i][1]
}

# Excel's ColumnWidth setter works in character-width units with ~1/7
# granularity, so 31.88 itself cannot be represented exactly; 31 is the
# closest input that lands nearest the source file's stored width (31.88).
$consolidated.Columns.Item(1).ColumnWidth = 31

# ---------------------------------------------------------------------
# Sheet3: "Ruoxuan Li" (new worksheet, added after "Fanjie Gao")
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ruoxuan = $wb.Worksheets.Add($null, $lastSheet)
$ruoxuan.Name = "Ruoxuan Li"

$ruoxuanRows = @(
    @("user sotry 1", "time estimation "),
    @("a frontend to place the app", "3h"),
    @("OAuth to spotify", "3h"),
    @("store the token", "3h"),
    @("test api for spotify", "2h"),
    @("test api for imdb", "2h"),
    @("get user playlist genre", "1h"),
    @("query imdb with the genres", "1h"),
    @("algorithm to sort the movies", "2h"),
    @("database to store user and their movie list", "3h"),
    @("option to remove movie form list", "1h"),
    @("excluded movie list in database", "1h")
)

$fanjie.Range("A1:A12").Copy() | Out-Null
$ruoxuan.Range("A1:A12").PasteSpecial(-4122) | Out-Null
$fanjie.Range("A1:A12").Copy() | Out-Null
$ruoxuan.Range("D1:D12").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $ruoxuanRows.Count; $i++) {
    $r = $i + 1
    $ruoxuan.Cells.Item($r, 1).Value = $ruoxuanRows[$i][0]
    $ruoxuan.Cells.Item($r, 4).Value = $ruoxuanRows[$i][1]
}
